# Inventario.xlsx maintenance edit:
#  - rename the report sheet to the new date
#  - add a bar chart + a pie chart backed by the inventory table
#  - add a colour-scale conditional format on the "Precio Total" column
#  - add a data validation dropdown on the "Nombre" column
#  - turn the inventory range into a real Excel Table

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1) Rename the dated report sheet
$ws2.Name = "Reporte 2025-04-23"

# 2) Bar chart: quantities per product (categories = A2:A6, values = D2:D6)
$chartObj1 = $ws1.ChartObjects().Add(400, 50, 425.2, 212.6)
$chart1 = $chartObj1.Chart
$chart1.ChartType = 51
$chart1.ChartStyle = 10
$chart1.SetSourceData($ws1.Range("D1:D6"))
$series1 = $chart1.SeriesCollection().Item(1)
$series1.XValues = "='Inventario principal'!`$A`$2:`$A`$6"
$series1.Values = "='Inventario principal'!`$D`$2:`$D`$6"
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "Cantidad de productos"

$catAxis1 = $chart1.Axes(1)
$catAxis1.HasTitle = $true
$catAxis1.AxisTitle.Text = "Productos"

$valAxis1 = $chart1.Axes(2)
$valAxis1.HasTitle = $true
$valAxis1.AxisTitle.Text = "Cantidad de productos"

# 3) Pie chart: category distribution
$chartObj2 = $ws1.ChartObjects().Add(400, 350, 425.2, 212.6)
$chart2 = $chartObj2.Chart
$chart2.ChartType = 5
$chart2.ChartStyle = 10
$chart2.SetSourceData($ws1.Range("D2:D6"))
$series2 = $chart2.SeriesCollection().Item(1)
$series2.XValues = "='Inventario principal'!`$A`$2:`$A`$6"
$series2.Values = "='Inventario principal'!`$D`$2:`$D`$6"
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Distribución de categorías"

# 4) Conditional formatting: colour scale on Precio Total
$cfRange = $ws1.Range("F2:F100")
$colorScale = $cfRange.FormatConditions.AddColorScale(3)

# 5) Data validation: dropdown list on Nombre
$dvRange = $ws1.Range("B2:B6")
$dvRange.Validation.Add(3, 1, 1, '"Producto A,Producto B,Producto C,Producto D,Producto E"')
$dvRange.Validation.IgnoreBlank = $false
$dvRange.Validation.InCellDropdown = $true
$dvRange.Validation.ShowInput = $false
$dvRange.Validation.ShowError = $false
$dvRange.Validation.ErrorTitle = "Error de validación"
$dvRange.Validation.ErrorMessage = "Por favor ingresar un producto correcto"

# 6) Turn the inventory range into a real Excel Table
$tableRange = $ws1.Range("A1:F6")
$table = $ws1.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "TablaInventario"
$table.TableStyle = "TableStyleMedium9"
$table.ShowTableStyleColumnStripes = $true

Write-Host "Inventario.xlsx updated"
